$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: add the evaluation result that used to sit alone in D18,
# plus a new "evaluation" label next to it.
$ws.Range("F3").Value = -0.048907
$ws.Range("G3").Value = "evaluation"

# The old stray value in D18 is obsolete now that it lives in F3.
$ws.Range("D18").ClearContents()

# Row 14 no longer carries the "best" label -- that moves down to the
# new best-policy row (17) below.
$ws.Range("D14").ClearContents()

# New row 16: a policy evaluation that "sucks".
$ws.Range("A16").Value = 0.006
$ws.Range("B16").Value = 0.002
$ws.Range("C16").Value = "sucks"

# New row 17: our best policy.
$ws.Range("A17").Value = 0.052
$ws.Range("B17").Value = 0.022
$ws.Range("C17").Value = 0.055397
$ws.Range("D17").Value = "best"

# Update the selection to reflect where the author left off.
$ws.Range("D5").Select()
